$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1232
$ws1.Range("F3").Value = 434
$ws1.Range("F5").Value = 12521
$ws1.Range("F6").Value = 69
$ws1.Range("F10").Value = 12393
$ws1.Range("F11").Value = 237
$ws1.Range("F12").Value = 4894
$ws1.Range("F13").Value = 4814
$ws1.Range("F14").Value = 151
$ws1.Range("F17").Value = 105
$ws1.Range("F18").Value = 962
$ws1.Range("F19").Value = 10

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1232
$ws4.Range("F3").Value = 434
$ws4.Range("F7").Value = 12521
$ws4.Range("F8").Value = 69
$ws4.Range("F12").Value = 12393
$ws4.Range("F13").Value = 237
$ws4.Range("F14").Value = 4894
$ws4.Range("F15").Value = 4814
$ws4.Range("F16").Value = 151
$ws4.Range("F19").Value = 105
$ws4.Range("F20").Value = 962
$ws4.Range("F21").Value = 10
